$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the "Date" value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# --- Elements sheet: swap the two "Mapping" columns (AK <-> AL) ---
$els = $wb.Worksheets.Item("Elements")

# Header row: the column titles swap places.
$els.Range("AK1").Value = "Mapping: Spécification métier vers l'extension ROR LocationStatus"
$els.Range("AL1").Value = "Mapping: RIM Mapping"

# Column widths swap along with the content (col 37 <-> col 38).
$els.Columns.Item(37).ColumnWidth = 70.56640625
$els.Columns.Item(38).ColumnWidth = 24.98046875

# Data rows: swap the cell contents between AK and AL for each row that had data.
# (A bare "'" forces an explicit empty-text cell instead of clearing it outright,
# matching the original empty shared string that lived in either column.)
$rows = @(3, 5, 6)
foreach ($r in $rows) {
    $akCell = $els.Range("AK$r")
    $alCell = $els.Range("AL$r")
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2

    if ($alVal -eq "") {
        $akCell.Value = "'"
    } else {
        $akCell.Value = $alVal
    }

    if ($akVal -eq "") {
        $alCell.Value = "'"
    } else {
        $alCell.Value = $akVal
    }
}
